$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 362, shifting existing rows 362-388 down to 363-389
$ws.Rows.Item(362).Insert()

# Populate the newly inserted row 362 with the new weekly record
$ws.Cells.Item(362, 1).Value = 9
$ws.Cells.Item(362, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(362, 3).Value = "Metropolitana"
$ws.Cells.Item(362, 4).Value = 44706
$ws.Cells.Item(362, 5).Value = 13
$ws.Cells.Item(362, 6).Value = 100112013
$ws.Cells.Item(362, 7).Value = "Alcachofa"
$ws.Cells.Item(362, 8).Value = "Española"
$ws.Cells.Item(362, 9).Value = "Primera"
$ws.Cells.Item(362, 10).Value = 110
$ws.Cells.Item(362, 11).Value = 24000
$ws.Cells.Item(362, 12).Value = 24000
$ws.Cells.Item(362, 13).Value = 24000
$ws.Cells.Item(362, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(362, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(362, 16).Value = 800
$ws.Cells.Item(362, 17).Value = 30
$ws.Cells.Item(362, 18).Value = "Hortaliza"
